$d = $word.ActiveDocument

# Locate the anchor paragraph (the last paragraph discussing VHF reflections / grounding)
# and collapse to its end so we can append a brand-new paragraph right after it.
$anchor = $d.Content.Duplicate()
[void]$anchor.Find.Execute("instead grounding through the connector is desired.", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
[void]$anchor.InsertParagraphAfter()

# The freshly inserted paragraph is now the last paragraph in the story; it already
# inherited the correct pPr/rPr (rFonts minorHAnsi theme + szCs 24) from its predecessor.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$fullText = "The design mistake that caused difficulty in tuning and grounding is that it was assumed that since the Oresat design did not require L-Matching to the 75-ohm quarter-wavelength monopole elements, our design would greatly improve with the inclusion of a simple matching network to step the 50-ohm output of the circuit to the 75-ohm impedance of the elements. This would most likely merge the 1st and 2nd resonance regions together such that it would be easier to match the antenna to 436 MHz."
$newRange.Text = $fullText

# Superscript the "st" in "1st"
$ordinal1 = $newPara.Range.Duplicate()
[void]$ordinal1.Find.Execute("1st", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$stRange = $d.Range($ordinal1.Start + 1, $ordinal1.End)
$stRange.Font.Superscript = $true

# Superscript the "nd" in "2nd"
$ordinal2 = $newPara.Range.Duplicate()
[void]$ordinal2.Find.Execute("2nd", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$ndRange = $d.Range($ordinal2.Start + 1, $ordinal2.End)
$ndRange.Font.Superscript = $true

Write-Output $d.Paragraphs.Last.Range.Text
